$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear previous data rows (2-6) fully so any removed cells are blanked,
# then rewrite rows 2-8 with the updated dataset.
$ws.Range("A2:AQ6").Clear()

# Row 2
$ws.Range("A2").Value = "Indonesia"
$ws.Range("B2").Value = "'6"
$ws.Range("C2").Value = "Insurance (Life)"
$ws.Range("D2").Value = -0.132
$ws.Range("E2").Value = 0.1278
$ws.Range("G2").Value = 0.134840857051875
$ws.Range("H2").Value = 0.134840857051875
$ws.Range("I2").Value = 0.127006535551276
$ws.Range("J2").Value = 0.1223171996968963
$ws.Range("K2").Value = 226.157
$ws.Range("L2").Value = 0.1714687552125191
$ws.Range("M2").Value = 10.5
$ws.Range("N2").Value = 0.003223627655655164
$ws.Range("O2").Value = 0.04642792396432567
$ws.Range("P2").Value = 10.5
$ws.Range("Q2").Value = 0.003223627655655164
$ws.Range("R2").Value = 0.04642792396432567
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 930.974
$ws.Range("V2").Value = 0.2858203364853248
$ws.Range("W2").Value = 0.04892263224368968
$ws.Range("X2").Value = 0.07261801482476032
$ws.Range("Y2").Value = -0.02369538258107064
$ws.Range("Z2").Value = 0.4694104164738877
$ws.Range("AA2").Value = 0.03392243662499179
$ws.Range("AB2").Value = 0.07261801482476032
$ws.Range("AC2").Value = -0.03870093540354935
$ws.Range("AD2").Value = 132.757
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 132.757
$ws.Range("AG2").Value = -798.2170000000001
$ws.Range("AH2").Value = 0.0391618536754301
$ws.Range("AI2").Value = 0.02719501350665664
$ws.Range("AJ2").Value = -0.3246126549065204
$ws.Range("AK2").Value = -0.2020448058100187
$ws.Range("AL2").Value = 73.303
$ws.Range("AM2").Value = 73.303
$ws.Range("AN2").Value = 0.7741655198152597
$ws.Range("AO2").Value = 2.285227071197632
$ws.Range("AP2").Value = -4.654760794009937
$ws.Range("AQ2").Value = 2.285227071197632

# Row 3
$ws.Range("A3").Value = "Indonesia"
$ws.Range("B3").Value = "PT Bhakti Multi Artha Tbk (IDX:BHAT)"
$ws.Range("C3").Value = "Insurance (Life)"
$ws.Range("G3").Value = 0.03832752613240418
$ws.Range("H3").Value = 0.03832752613240418
$ws.Range("I3").Value = 0.04912891986062717
$ws.Range("J3").Value = 0.04912891986062717
$ws.Range("K3").Value = 0.281
$ws.Range("L3").Value = 0.04895470383275262
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = -0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("S3").Value = 0
$ws.Range("U3").Value = 0.734
$ws.Range("V3").Value = 0.003288530465949821
$ws.Range("X3").Value = 0.07261801482476032
$ws.Range("AB3").Value = 0.07261801482476032
$ws.Range("AD3").Value = 0
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = -0.734
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = -0.003299380579504284
$ws.Range("AK3").Value = -0.02135831926904498
$ws.Range("AL3").Value = 0.003
$ws.Range("AM3").Value = 0.003
$ws.Range("AN3").Value = 0
$ws.Range("AO3").Value = 93.99999999999999
$ws.Range("AP3").Value = -2.430463576158941
$ws.Range("AQ3").Value = 93.99999999999999

# Row 4
$ws.Range("A4").Value = "Indonesia"
$ws.Range("B4").Value = "PT Capital Financial Indonesia Tbk (IDX:CASA)"
$ws.Range("C4").Value = "Insurance (Life)"
$ws.Range("G4").Value = 0.06495922163399628
$ws.Range("H4").Value = 0.06495922163399628
$ws.Range("I4").Value = 0.1070253255115181
$ws.Range("J4").Value = 0.1070253255115181
$ws.Range("K4").Value = 3.44
$ws.Range("L4").Value = 0.004922020317642009
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = -0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = -0
$ws.Range("S4").Value = 0
$ws.Range("U4").Value = 141.8
$ws.Range("V4").Value = 0.09614863032275564
$ws.Range("W4").Value = 0.00786825251601098
$ws.Range("X4").Value = 0.07743877919592808
$ws.Range("Y4").Value = -0.0695705266799171
$ws.Range("Z4").Value = 4.615333817605494
$ws.Range("AA4").Value = 0.4939576041735455
$ws.Range("AB4").Value = 0.07405098413322345
$ws.Range("AC4").Value = 0.4199066200403221
$ws.Range("AD4").Value = 132.1
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 132.1
$ws.Range("AG4").Value = -9.700000000000017
$ws.Range("AH4").Value = 0.08220797809446762
$ws.Range("AI4").Value = 0.1993962264150943
$ws.Range("AJ4").Value = -0.006620708484062533
$ws.Range("AK4").Value = -0.01862876896485504
$ws.Range("AL4").Value = 73.3
$ws.Range("AM4").Value = 73.3
$ws.Range("AN4").Value = 1.763684913217623
$ws.Range("AO4").Value = 1.020463847203274
$ws.Range("AP4").Value = -0.1295060080106811
$ws.Range("AQ4").Value = 1.020463847203274

# Row 5
$ws.Range("A5").Value = "Indonesia"
$ws.Range("B5").Value = "PT Paninvest Tbk (IDX:PNIN)"
$ws.Range("C5").Value = "Insurance (Life)"
$ws.Range("D5").Value = -0.156
$ws.Range("E5").Value = 0.0946
$ws.Range("G5").Value = 0.2914653784219001
$ws.Range("H5").Value = 0.2914653784219001
$ws.Range("I5").Value = 0.1604938271604938
$ws.Range("J5").Value = 0.1601440943430993
$ws.Range("K5").Value = 72.5
$ws.Range("L5").Value = 0.3891572732152442
$ws.Range("M5").Value = -0
$ws.Range("N5").Value = -0
$ws.Range("O5").Value = -0
$ws.Range("P5").Value = -0
$ws.Range("Q5").Value = -0
$ws.Range("R5").Value = -0
$ws.Range("S5").Value = 0
$ws.Range("U5").Value = 374.7
$ws.Range("V5").Value = 1.478689818468824
$ws.Range("W5").Value = 0.06349623401646523
$ws.Range("X5").Value = 0.07261801482476032
$ws.Range("Y5").Value = -0.009121780808295094
$ws.Range("Z5").Value = 0.2255174918290764
$ws.Range("AA5").Value = 0.03611529448749474
$ws.Range("AB5").Value = 0.07261801482476032
$ws.Range("AC5").Value = -0.03650272033726559
$ws.Range("AD5").Value = 0
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 0
$ws.Range("AG5").Value = -374.7
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 3.089035449299258
$ws.Range("AK5").Value = -0.2476209357652656
$ws.Range("AL5").Value = 0
$ws.Range("AM5").Value = 0
$ws.Range("AN5").Value = 0
$ws.Range("AP5").Value = -12.36633663366337

# Row 6
$ws.Range("A6").Value = "Indonesia"
$ws.Range("B6").Value = "PT Asuransi Jiwa Sinarmas MSIG Tbk (IDX:LIFE)"
$ws.Range("C6").Value = "Insurance (Life)"
$ws.Range("G6").Value = 0.1195652173913044
$ws.Range("H6").Value = 0.1195652173913044
$ws.Range("I6").Value = 0.08260869565217391
$ws.Range("J6").Value = 0.06754261635701972
$ws.Range("K6").Value = 18.6
$ws.Range("L6").Value = 0.0808695652173913
$ws.Range("M6").Value = 10.5
$ws.Range("N6").Value = 0.01449875724937862
$ws.Range("O6").Value = 0.564516129032258
$ws.Range("P6").Value = 10.5
$ws.Range("Q6").Value = 0.01449875724937862
$ws.Range("R6").Value = 0.564516129032258
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 64.2
$ws.Range("V6").Value = 0.08864954432477216
$ws.Range("W6").Value = 0.03434903047091413
$ws.Range("X6").Value = 0.07266684104041744
$ws.Range("Y6").Value = -0.03831781056950331
$ws.Range("Z6").Value = 0.4697712418300654
$ws.Range("AA6").Value = 0.03172957876248884
$ws.Range("AB6").Value = 0.07262872923232194
$ws.Range("AC6").Value = -0.04089915046983311
$ws.Range("AD6").Value = 0.657
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 0.657
$ws.Range("AG6").Value = -63.54300000000001
$ws.Range("AH6").Value = 0.0009063856733121153
$ws.Range("AI6").Value = 0.001262594718625866
$ws.Range("AJ6").Value = -0.09618152838765048
$ws.Range("AK6").Value = -0.1393007232159103
$ws.Range("AL6").Value = 0
$ws.Range("AM6").Value = 0
$ws.Range("AN6").Value = 0.03173913043478261
$ws.Range("AP6").Value = -3.069710144927537

# Row 7
$ws.Range("A7").Value = "Indonesia"
$ws.Range("B7").Value = "PT Panin Financial Tbk (IDX:PNLF)"
$ws.Range("C7").Value = "Insurance (Life)"
$ws.Range("D7").Value = -0.108
$ws.Range("E7").Value = 0.161
$ws.Range("G7").Value = 0.2574206755373593
$ws.Range("H7").Value = 0.2574206755373593
$ws.Range("I7").Value = 0.2221084953940634
$ws.Range("J7").Value = 0.2155174709372477
$ws.Range("K7").Value = 131.2
$ws.Range("L7").Value = 0.6714431934493347
$ws.Range("M7").Value = -0
$ws.Range("N7").Value = -0
$ws.Range("O7").Value = -0
$ws.Range("P7").Value = -0
$ws.Range("Q7").Value = -0
$ws.Range("R7").Value = -0
$ws.Range("S7").Value = 0
$ws.Range("U7").Value = 347.3
$ws.Range("V7").Value = 0.6123060648801129
$ws.Range("W7").Value = 0.07997561718988112
$ws.Range("X7").Value = 0.07261801482476032
$ws.Range("Y7").Value = 0.007357602365120799
$ws.Range("Z7").Value = 0.1460934579439252
$ws.Range("AA7").Value = 0.03148569257655193
$ws.Range("AB7").Value = 0.07261801482476032
$ws.Range("AC7").Value = -0.04113232224820839
$ws.Range("AD7").Value = 0
$ws.Range("AE7").Value = 0
$ws.Range("AF7").Value = 0
$ws.Range("AG7").Value = -347.3
$ws.Range("AH7").Value = 0
$ws.Range("AI7").Value = 0
$ws.Range("AJ7").Value = -1.579354251932696
$ws.Range("AK7").Value = -0.2444225490886058
$ws.Range("AL7").Value = 0
$ws.Range("AM7").Value = 0
$ws.Range("AN7").Value = 0
$ws.Range("AP7").Value = -7.700665188470066

# Row 8
$ws.Range("A8").Value = "Indonesia"
$ws.Range("B8").Value = "PT Asuransi Jiwa Syariah Jasa Mitra Abadi Tbk (IDX:JMAS)"
$ws.Range("C8").Value = "Insurance (Life)"
$ws.Range("G8").Value = 0.04884615384615384
$ws.Range("H8").Value = 0.04884615384615384
$ws.Range("I8").Value = 0.05076923076923077
$ws.Range("J8").Value = 0.05039865244244806
$ws.Range("K8").Value = 0.136
$ws.Range("L8").Value = 0.05230769230769231
$ws.Range("M8").Value = -0
$ws.Range("N8").Value = -0
$ws.Range("O8").Value = -0
$ws.Range("P8").Value = -0
$ws.Range("Q8").Value = -0
$ws.Range("R8").Value = -0
$ws.Range("S8").Value = 0
$ws.Range("U8").Value = 2.24
$ws.Range("V8").Value = 0.1555555555555556
$ws.Range("W8").Value = 0.0171500630517024
$ws.Range("X8").Value = 0.07261801482476032
$ws.Range("Y8").Value = -0.05546795177305792
$ws.Range("Z8").Value = 0.5048543689320388
$ws.Range("AA8").Value = 0.02544397987385727
$ws.Range("AB8").Value = 0.07261801482476032
$ws.Range("AC8").Value = -0.04717403495090305
$ws.Range("AD8").Value = 0
$ws.Range("AE8").Value = 0
$ws.Range("AF8").Value = 0
$ws.Range("AG8").Value = -2.24
$ws.Range("AH8").Value = 0
$ws.Range("AI8").Value = 0
$ws.Range("AJ8").Value = -0.1842105263157895
$ws.Range("AK8").Value = -0.4171322160148976
$ws.Range("AL8").Value = 0
$ws.Range("AM8").Value = 0
$ws.Range("AN8").Value = 0
$ws.Range("AP8").Value = -12.30769230769231
